$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("VEDA_Sets-Proc")

# Update the PSET_PN pattern for the "Util Batt Stg" set (row 19)
$ws.Range("B19").Value = "EN*STG?hb*,-*EV*"

# Add the And/Or qualifiers for row 19, matching the pattern used by other
# multi-token PSET_PN rows (rows 3, 7, 17)
$ws.Range("H19").Value = "And"
$ws.Range("I19").Value = "Or"

# Append a new PSET row (21) defining the "Grid" set
$ws.Range("A21").Value = "IRE"
$ws.Range("B21").Value = "g[_]*"
$ws.Range("F21").Value = "Grid"

$wb.Save()
